# Update 2p0. Convention change to support multi-axle vehicles.
# Rename the axle labels sAxleF/sAxleR -> sAxle1/sAxle2 on every
# Body sheet so additional axles beyond front/rear can be added later.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Sedan_HambaLG", "Sedan_Hamba", "Bus_Makhulu")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A5").Value = "sAxle1"
    $ws.Range("A6").Value = "sAxle2"
}

# Update view state: Sedan_HambaLG becomes the active/selected tab
# (previously Bus_Makhulu was active), and update each sheet's
# current selection.
$wsHamba = $wb.Worksheets.Item("Sedan_Hamba")
$wsHamba.Range("A17").Select()

$wsMakhulu = $wb.Worksheets.Item("Bus_Makhulu")
$wsMakhulu.Range("A2").Select()

$wsHambaLG = $wb.Worksheets.Item("Sedan_HambaLG")
$wsHambaLG.Activate()
$wsHambaLG.Range("A2").Select()
